# Insert a new weekly price-report row for "Femacal de La Calera - Cilantro"
# at row 403, pushing the existing rows (403..429) down to (404..430).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 403 (existing row 403 and below shift down by one).
$ws.Rows.Item(403).Insert()

# Populate the newly inserted row 403 with the new observation.
$ws.Cells.Item(403, 1).Value = 3
$ws.Cells.Item(403, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(403, 3).Value = "Coquimbo"
$ws.Cells.Item(403, 4).Value = 44826
$ws.Cells.Item(403, 5).Value = 5
$ws.Cells.Item(403, 6).Value = 100112040
$ws.Cells.Item(403, 7).Value = "Cilantro"
$ws.Cells.Item(403, 8).Value = "Sin especificar"
$ws.Cells.Item(403, 9).Value = "Primera"
$ws.Cells.Item(403, 10).Value = 190
$ws.Cells.Item(403, 11).Value = 4000
$ws.Cells.Item(403, 12).Value = 4500
$ws.Cells.Item(403, 13).Value = 4211
$ws.Cells.Item(403, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(403, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(403, 16).Value = 1404
$ws.Cells.Item(403, 17).Value = 3
$ws.Cells.Item(403, 18).Value = "Hortaliza"
